$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.817.74"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.557.74"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.94"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.26"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.53%  "
$ws.Range("E7").Value = "  -0.78%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.56"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.951.04"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("E15").Value = "  +5.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.604.90"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.838"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.815.90"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.74"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0958"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.35"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.50"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "247.51"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.92"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.79"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.40"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.98"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("E30").Value = "  -1.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.15"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.75"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0795"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.08"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.94%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.72"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.64"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.58"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +11.13%  "
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.69"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.06"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.53%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.988.13"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.98%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.20"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.96"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.804.12"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "81.37"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.193"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.41"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.47%  "
